$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Files")

# Insert a new column before column E (shifts Subject..Label from E:P to F:Q)
$ws.Range("E:E").Insert()

# Set header for the newly inserted column
$ws.Range("E1").Value = "SessionLabel"

# Update selection to match the target view (E2, no frozen topLeftCell)
$ws.Range("E2").Select()
